$d = $word.ActiveDocument

$replacements = @(
    @('2025-09-28 Sunday', '2025-09-29 Monday'),
    @('2+27=', '69-43='),
    @('51-1=', '33+53='),
    @('97-34=', '17+53='),
    @('96-92=', '60+36='),
    @('46-30=', '18+32='),
    @('77-34=', '36+63='),
    @('4+34=', '64+29='),
    @('69-35=', '37-15='),
    @('64+24=', '97-32='),
    @('46-4=', '98-18='),
    @('36+35=', '26+3='),
    @('78-72=', '37-30='),
    @('3-0=', '5+88='),
    @('45+23=', '59-15='),
    @('11+40=', '26-15='),
    @('43+25=', '97+0='),
    @('60+26=', '3+93='),
    @('26-7=', '6+31='),
    @('74-0=', '77-49='),
    @('49-47=', '11+19='),
    @('16+66=', '15+18='),
    @('46+25=', '2-0='),
    @('2+88=', '53+33='),
    @('69-47=', '8+20='),
    @('27+0=', '76+5='),
    @('8+35=', '50+26='),
    @('90-4=', '57+2='),
    @('73-69=', '46+5='),
    @('78+13=', '33+46='),
    @('48+3=', '0+8='),
    @('4+39=', '94-82='),
    @('64-11=', '42-3='),
    @('78-18=', '58+19='),
    @('40+22=', '37-13='),
    @('16+34=', '19+80='),
    @('20-13=', '32-18='),
    @('37+33=', '46+5='),
    @('91-29=', '81-40='),
    @('87-84=', '30+67='),
    @('93-42=', '99-55='),
    @('14+10=', '69+21='),
    @('30-16=', '56+37='),
    @('70-40=', '49+22='),
    @('46-3=', '96-5='),
    @('79-33=', '37-8='),
    @('73-58=', '12+69='),
    @('17+13=', '6+2='),
    @('47+45=', '53+18='),
    @('84-18=', '35+16='),
    @('15+52=', '26+66='),
    @('15-12=', '59-43='),
    @('44+35=', '85-32='),
    @('11+31=', '5+93='),
    @('41+42=', '75+12='),
    @('64+32=', '15+47='),
    @('43+42=', '97-85='),
    @('25+33=', '60-36='),
    @('20+57=', '1-0='),
    @('49+28=', '11+14='),
    @('36-25=', '2+34='),
    @('79-38=', '39-5='),
    @('38+32=', '72-1='),
    @('17+76=', '30-27='),
    @('10+19=', '66-14='),
    @('53+6=', '26+37='),
    @('42+50=', '92-18='),
    @('79-8=', '95-63='),
    @('28+3=', '46-10='),
    @('66+16=', '21+34='),
    @('75-30=', '21+32='),
    @('14+74=', '67+2='),
    @('18+60=', '74-47='),
    @('77-15=', '99-56='),
    @('88-12=', '2+71='),
    @('68-22=', '22-20='),
    @('46+44=', '32+64='),
    @('40+4=', '29-22='),
    @('29+40=', '87-51='),
    @('19+31=', '8+41='),
    @('67+6=', '93-36='),
    @('2+9=', '37+10='),
    @('42+40=', '26+32='),
    @('3+12=', '33-20='),
    @('38+5=', '91-91='),
    @('28+40=', '79-29='),
    @('84-82=', '95-30='),
    @('75-34=', '17+11='),
    @('20+2=', '82+4='),
    @('3+17=', '16+63='),
    @('65-9=', '52-48='),
    @('54-36=', '12+70='),
    @('18-5=', '49+10='),
    @('54-32=', '17+16='),
    @('36-6=', '45+32='),
    @('65-21=', '92-77='),
    @('87-57=', '65+14='),
    @('21+28=', '91+8='),
    @('5+15=', '56+1='),
    @('34+8=', '54+27='),
    @('10+28=', '11+53='),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

"Done: applied $($replacements.Count) replacements"